$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.545999999999999
$ws.Range("A3").Value = -21.937
$ws.Range("E3").Value = 16.483
$ws.Range("E12").Value = 17.692
$ws.Range("A14").Value = -21.659
$ws.Range("A21").Value = -19.861
$ws.Range("A23").Value = -20.317
$ws.Range("E24").Value = 16.923
$ws.Range("A25").Value = -21.608
$ws.Range("D25").Value = -7.377999999999998
$ws.Range("E25").Value = 17.079
$ws.Range("A26").Value = -21.006
$ws.Range("D27").Value = -8.372
$ws.Range("A29").Value = -21.297
$ws.Range("D31").Value = -8.222
$ws.Range("D39").Value = -7.502
$ws.Range("D48").Value = -7.475
$ws.Range("E50").Value = 16.45
$ws.Range("D51").Value = -8.294
$ws.Range("D52").Value = -8.083
$ws.Range("A53").Value = -21.945
$ws.Range("E53").Value = 16.484
$ws.Range("D55").Value = -8.035
$ws.Range("D56").Value = -8.409000000000001
$ws.Range("A57").Value = -22.563
$ws.Range("D57").Value = -8.106999999999999
$ws.Range("E57").Value = 16.481
$ws.Range("A59").Value = -22.5
$ws.Range("E61").Value = 16.602
$ws.Range("E63").Value = 17.852
$ws.Range("A69").Value = -21.626
$ws.Range("E70").Value = 17.547
$ws.Range("D73").Value = -8.004000000000001
$ws.Range("A79").Value = -20.856
$ws.Range("A83").Value = -21.997
$ws.Range("E86").Value = 16.394
$ws.Range("D89").Value = -6.849000000000001
$ws.Range("D90").Value = -7.475999999999999
$ws.Range("A91").Value = -21.493
$ws.Range("D92").Value = -6.794
$ws.Range("A93").Value = -21.457
$ws.Range("E98").Value = 16.494
$ws.Range("E100").Value = 16.617
$ws.Range("E102").Value = 16.49
